# Auto-generated Excel COM-interop script to apply market-price refresh
# to the Spriggan_Profits workbook (columns H-N across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 25099.5
$ws.Range("J26").Value = 48999
$ws.Range("L26").Value = 48999
$ws.Range("N26").Value = -49687
$ws.Range("H28").Value = 997.64703
$ws.Range("I28").Value = 791.7
$ws.Range("K28").Value = 791.7
$ws.Range("M28").Value = -306.7
$ws.Range("H40").Value = 3249.2917
$ws.Range("I40").Value = 2598.4666
$ws.Range("K40").Value = 2598.4666
$ws.Range("M40").Value = -2423.4666
$ws.Range("H41").Value = 1736.4706
$ws.Range("J41").Value = 1616.6666
$ws.Range("L41").Value = 1616.6666
$ws.Range("N41").Value = -2496.6666
$ws.Range("H53").Value = 402
$ws.Range("I53").Value = 359.85715
$ws.Range("K53").Value = 359.85715
$ws.Range("M53").Value = 277.14285
$ws.Range("H141").Value = 3539.6
$ws.Range("I141").Value = 2924.75
$ws.Range("K141").Value = 8774.25
$ws.Range("M141").Value = -3594.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 507547.6
$ws.Range("I2").Value = 774301.8
$ws.Range("J2").Value = 714.5
$ws.Range("K2").Value = 774301.8
$ws.Range("L2").Value = 714.5
$ws.Range("M2").Value = -774188.8
$ws.Range("N2").Value = -940.5
$ws.Range("H46").Value = 22997.25
$ws.Range("J46").Value = 22998
$ws.Range("L46").Value = 22998
$ws.Range("N46").Value = -23636
$ws.Range("H116").Value = 507547.6
$ws.Range("I116").Value = 774301.8
$ws.Range("J116").Value = 714.5
$ws.Range("K116").Value = 774301.8
$ws.Range("L116").Value = 714.5
$ws.Range("M116").Value = -772007.8
$ws.Range("N116").Value = -5302.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 507547.6
$ws.Range("I3").Value = 774301.8
$ws.Range("J3").Value = 714.5
$ws.Range("K3").Value = 774301.8
$ws.Range("L3").Value = 714.5
$ws.Range("M3").Value = -774187.8
$ws.Range("N3").Value = -942.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9301
$ws.Range("I31").Value = 6351.2915
$ws.Range("K31").Value = 6351.2915
$ws.Range("M31").Value = -6056.2915
$ws.Range("H34").Value = 9301
$ws.Range("I34").Value = 6351.2915
$ws.Range("K34").Value = 6351.2915
$ws.Range("M34").Value = -6149.2915
$ws.Range("H99").Value = 1712
$ws.Range("I99").Value = 1528.2858
$ws.Range("K99").Value = 1528.2858
$ws.Range("M99").Value = -30.28580000000011
$ws.Range("H126").Value = 1712
$ws.Range("I126").Value = 1528.2858
$ws.Range("K126").Value = 4584.857400000001
$ws.Range("M126").Value = -2114.857400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 137768
$ws.Range("J37").Value = 137768
$ws.Range("L37").Value = 413304
$ws.Range("N37").Value = -413528
$ws.Range("H76").Value = 10000
$ws.Range("I76").Value = 10000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 30000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -29617
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 10000
$ws.Range("I79").Value = 10000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 30000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -28674
$ws.Range("N79").ClearContents()
$ws.Range("H129").Value = 1529.4615
$ws.Range("I129").Value = 609.3333
$ws.Range("J129").Value = 3599.75
$ws.Range("K129").Value = 1827.9999
$ws.Range("L129").Value = 10799.25
$ws.Range("M129").Value = 3172.0001
$ws.Range("N129").Value = -20799.25
$ws.Range("H131").Value = 1338.909
$ws.Range("I131").Value = 1153.625
$ws.Range("J131").Value = 1833
$ws.Range("K131").Value = 3460.875
$ws.Range("L131").Value = 5499
$ws.Range("M131").Value = 1579.125
$ws.Range("N131").Value = -15579

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 831.7143
$ws.Range("I2").Value = 1514
$ws.Range("J2").Value = 149.42857
$ws.Range("K2").Value = 1514
$ws.Range("L2").Value = 149.42857
$ws.Range("M2").Value = -1401
$ws.Range("N2").Value = -375.42857
$ws.Range("H70").Value = 7097.5386
$ws.Range("J70").Value = 5274.3335
$ws.Range("L70").Value = 5274.3335
$ws.Range("N70").Value = -5814.3335
$ws.Range("H73").Value = 7097.5386
$ws.Range("J73").Value = 5274.3335
$ws.Range("L73").Value = 5274.3335
$ws.Range("N73").Value = -7146.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1822.5714
$ws.Range("I40").Value = 1822.5714
$ws.Range("K40").Value = 1822.5714
$ws.Range("M40").Value = -1686.5714
$ws.Range("H46").Value = 1064.8572
$ws.Range("I46").Value = 939
$ws.Range("J46").Value = 1379.5
$ws.Range("K46").Value = 939
$ws.Range("L46").Value = 1379.5
$ws.Range("M46").Value = -751
$ws.Range("N46").Value = -1755.5
$ws.Range("H68").Value = 4168766.2
$ws.Range("J68").Value = 2649.5
$ws.Range("L68").Value = 2649.5
$ws.Range("N68").Value = -4147.5
$ws.Range("H71").Value = 4168766.2
$ws.Range("J71").Value = 2649.5
$ws.Range("L71").Value = 13247.5
$ws.Range("N71").Value = -20735.5
$ws.Range("H100").Value = 7679710.5
$ws.Range("I100").Value = 11742323
$ws.Range("J100").Value = 5887.3335
$ws.Range("K100").Value = 11742323
$ws.Range("L100").Value = 5887.3335
$ws.Range("M100").Value = -11741782
$ws.Range("N100").Value = -6969.3335
$ws.Range("H122").Value = 4571.533
$ws.Range("I122").Value = 4449.7036
$ws.Range("K122").Value = 13349.1108
$ws.Range("M122").Value = -10899.1108
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3447.4285
$ws.Range("I81").Value = 3763.3333
$ws.Range("J81").Value = 1552
$ws.Range("K81").Value = 7526.6666
$ws.Range("L81").Value = 3104
$ws.Range("M81").Value = -6465.6666
$ws.Range("N81").Value = -5226
$ws.Range("H84").Value = 3447.4285
$ws.Range("I84").Value = 3763.3333
$ws.Range("J84").Value = 1552
$ws.Range("K84").Value = 37633.333
$ws.Range("L84").Value = 15520
$ws.Range("M84").Value = -32329.333
$ws.Range("N84").Value = -26128
$ws.Range("H100").Value = 788.5714
$ws.Range("I100").Value = 896
$ws.Range("J100").Value = 520
$ws.Range("K100").Value = 1792
$ws.Range("L100").Value = 1040
$ws.Range("M100").Value = -1251
$ws.Range("N100").Value = -2122
$ws.Range("H113").Value = 564.3333
$ws.Range("I113").Value = 564.3333
$ws.Range("K113").Value = 1692.9999
$ws.Range("M113").Value = 477.0001
$ws.Range("H136").Value = 55560110
$ws.Range("I136").Value = 83337576
$ws.Range("K136").Value = 250012728
$ws.Range("M136").Value = -250010178
